# Apply "changes in budget and gdpr" commit:
# Increase the LOSSES figures for Year 1/2/3 by 15000 each.
# Dependent formulas (G4:H6) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E4").Value = 44968
$ws.Range("E5").Value = 28159
$ws.Range("E6").Value = 29023

# Update the active selection on the sheet (matches the saved view state).
$ws.Range("J6").Select()
